$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# New column G width (stored OOXML width 17 == ColumnWidth 17 - 5/6)
$ws.Columns.Item(7).ColumnWidth = 17 - 5/6

# Copy formatting from column F into the new column G so the new cells
# inherit the same styles used throughout the sheet (header / data / total).
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("F2:F54").Copy()
$ws.Range("G2:G54").PasteSpecial(-4122)

$ws.Range("F55").Copy()
$ws.Range("G55").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Header text
$ws.Cells.Item(1, 7).Value = "PRESUPUESTO"

# Data + total values (all zero per the source data)
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
}
